# Updated cryptos list on Fri Dec  1 20:06:44 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($range, $value)
    # Force the cell to keep/become plain text (matches the source workbook,
    # where Price/Volume columns are always stored as literal text strings),
    # then restore the default "Normal" style so no extra style index /
    # formatting is introduced on the cell.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "38.766.47"
Set-TextCell "E2" "  +2.60%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.090.04"
Set-TextCell "E3" "  +2.37%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  +0.02%  "

# Row 5 - BNB
Set-TextCell "D5" "228.82"
Set-TextCell "E5" "  +0.53%  "

# Row 6 - XRP
Set-TextCell "E6" "  +1.15%  "

# Row 7 - Solana
Set-TextCell "D7" "60.64"
Set-TextCell "E7" "  +0.82%  "

# Row 8 - USDC
Set-TextCell "E8" "  -0.02%  "

# Row 9 - Cardano
Set-TextCell "D9" "0.386"
Set-TextCell "E9" "  +2.56%  "

# Row 10 - Dogecoin
Set-TextCell "E10" "  -0.88%  "

# Row 11 - TRON
Set-TextCell "E11" "  -0.25%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextCell "D12" "2.399.82"
Set-TextCell "E12" "  +2.31%  "

# Row 13 - Chainlink
Set-TextCell "D13" "14.99"
Set-TextCell "E13" "  +4.37%  "

# Row 14 - Avalanche
Set-TextCell "D14" "21.87"
Set-TextCell "E14" "  +3.67%  "

# Row 15 - Polygon
Set-TextCell "E15" "  +4.45%  "

# Row 16 - Polkadot
Set-TextCell "E16" "  +0.15%  "

# Row 17 - WrappedEther
Set-TextCell "D17" "2.092.97"
Set-TextCell "E17" "  +1.90%  "

# Row 18 - WrappedBTC
Set-TextCell "D18" "38.711.91"
Set-TextCell "E18" "  +2.51%  "

# Row 19 - Litecoin
Set-TextCell "D19" "71.67"
Set-TextCell "E19" "  +3.26%  "

# Row 20 - Uniswap
Set-TextCell "D20" "6.05"
Set-TextCell "E20" "  +2.22%  "

# Row 21 - ShibaInu
Set-TextCell "D21" "0.0₃0837"
Set-TextCell "E21" "  +1.08%  "

# Row 22 - BitcoinCash
Set-TextCell "D22" "226.84"
Set-TextCell "E22" "  +1.40%  "

# Row 24 - Toncoin
Set-TextCell "E24" "  -0.25%  "

# Row 25 - PancakeSwap
Set-TextCell "E25" "  +3.03%  "

# Row 26 - Monero
Set-TextCell "E26" "  +1.20%  "

# Row 27 - Cosmos
Set-TextCell "D27" "9.45"
Set-TextCell "E27" "  +1.00%  "

# Row 28 - Kaspa
Set-TextCell "E28" "  +7.99%  "

# Row 29 - ImmutableX
Set-TextCell "D29" "1.44"
Set-TextCell "E29" "  +12.29%  "

# Row 30 - EthereumClassic
Set-TextCell "D30" "19.18"
Set-TextCell "E30" "  +2.20%  "

# Row 31 - Stellar
Set-TextCell "E31" "  +1.30%  "

# Row 32 - WEMIXToken
Set-TextCell "E32" "  +4.42%  "

# Row 33 - Filecoin
Set-TextCell "D33" "4.50"
Set-TextCell "E33" "  +2.95%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextCell "E34" "  +4.82%  "

# Row 35 - Hedera
Set-TextCell "E35" "  +2.18%  "

# Row 36 - was THORChain, now LidoDAOToken
Set-TextCell "B36" "LidoDAOToken"
Set-TextCell "C36" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D36" "2.40"
Set-TextCell "E36" "  +2.90%  "

# Row 37 - was LidoDAOToken, now THORChain
Set-TextCell "B37" "THORChain"
Set-TextCell "C37" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell "D37" "6.44"
Set-TextCell "E37" "  -2.97%  "

# Row 38 - RenderToken
Set-TextCell "D38" "3.56"
Set-TextCell "E38" "  +3.10%  "

# Row 39 - BinanceUSD
Set-TextCell "E39" "  -0.04%  "

# Row 40 - InjectiveProtocol
Set-TextCell "E40" "  +1.67%  "

# Row 41 - Maker
Set-TextCell "D41" "1.543.76"
Set-TextCell "E41" "  +0.66%  "

# Row 42 - Aave
Set-TextCell "D42" "101.04"
Set-TextCell "E42" "  +3.20%  "

# Row 43 - VeChain
Set-TextCell "D43" "0.0224"
Set-TextCell "E43" "  +4.28%  "

# Row 44 - was Cronos, now HuobiToken
Set-TextCell "B44" "HuobiToken"
Set-TextCell "C44" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell "D44" "2.82"
Set-TextCell "E44" "  -0.81%  "

# Row 45 - was HuobiToken, now Cronos
Set-TextCell "B45" "Cronos"
Set-TextCell "C45" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D45" "0.0922"
Set-TextCell "E45" "  +1.91%  "

# Row 46 - FraxShare
Set-TextCell "D46" "7.67"
Set-TextCell "E46" "  +9.41%  "

# Row 47 - TrustWalletToken
Set-TextCell "D47" "1.13"
Set-TextCell "E47" "  +1.48%  "

# Row 48 - FTXToken
Set-TextCell "E48" "  -0.79%  "

# Row 49 - ARBITRUM
Set-TextCell "E49" "  +2.69%  "

# Row 50 - MXToken
Set-TextCell "E50" "  +1.13%  "

# Row 51 - RocketPoolETH
Set-TextCell "D51" "2.288.61"
Set-TextCell "E51" "  +2.44%  "
